$wb = $excel.ActiveWorkbook

# Work on the "Part3-2" sheet (sheet3.xml)
$ws3 = $wb.Worksheets.Item("Part3-2")

# Add the new data point in row 3
$ws3.Range("E3").Value = 4283

# Select E3 and make this sheet the active tab
$ws3.Activate()
$ws3.Range("E3").Select() | Out-Null
